# Weekly update: insert a new "Black Amber / Primera" price record ahead of
# the existing history (new row 19) and push the rest of the table down one
# row, matching the upstream "Fruta / hortaliza, semanal" refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 19:42 down to 20:43 (Excel copies formatting from the row
# above automatically, same as a manual "Insert Copied Cells" / right-click
# "Insert" on the row header).
$ws.Rows("19:19").Insert()

# Populate the newly inserted row 19 with the new weekly record. Columns
# A,B,C,E,F,G,H,I,J,K,L carry over unchanged from the record that used to
# occupy row 19 (now row 20), so just re-assert them explicitly for clarity.
$ws.Range("A19").Value = 11
$ws.Range("B19").Value = "Vega Monumental Concepción"
$ws.Range("C19").Value = "Bíobío"
$ws.Range("D19").Value = 44580
$ws.Range("E19").Value = 8
$ws.Range("F19").Value = "Fruta"
$ws.Range("G19").Value = 100103
$ws.Range("H19").Value = "Frutos de hueso (carozo)"
$ws.Range("I19").Value = 100103002
$ws.Range("J19").Value = "Ciruela"
$ws.Range("K19").Value = "Black Amber"
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 230
$ws.Range("N19").Value = 10000
$ws.Range("O19").Value = 11000
$ws.Range("P19").Value = 10565
$ws.Range("Q19").Value = "$/bandeja 18 kilos granel"
$ws.Range("R19").Value = "Provincia de Curicó"
$ws.Range("S19").Value = 587
$ws.Range("T19").Value = 18
